# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 08:05"

# Row 42 - Israel
$ws.Range("B42").Value = 16771
$ws.Range("C42").Value = 14
$ws.Range("D42").Value = 14486
$ws.Range("E42").Value = 2004

# Row 48 - Afganistan
$ws.Range("B48").Value = 12456
$ws.Range("C48").Value = 625
$ws.Range("D48").Value = 1138
$ws.Range("E48").Value = 11091
$ws.Range("G48").Value = 7
$ws.Range("H48").Value = 227

# Row 59 - Oman
$ws.Range("E59").Value = 6013
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 38

# Row 84 - Bulgaria
$ws.Range("B84").Value = 2460
$ws.Range("C84").Value = 17
$ws.Range("D84").Value = 912
$ws.Range("E84").Value = 1415
$ws.Range("G84").Value = 3
$ws.Range("H84").Value = 133

# Row 89 - El Salvador
$ws.Range("B89").Value = 2109
$ws.Range("C89").Value = 67
$ws.Range("D89").Value = 873
$ws.Range("E89").Value = 1199
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 37

# Row 148 - Vietnam
$ws.Range("D148").Value = 278
$ws.Range("E148").Value = 49

# Rows 207/208 - Groenlandia and Islas Turcas y Caicos swap order/data
$ws.Range("A207").Value = "Groenlandia"
$ws.Range("D207").Value = 11
$ws.Range("H207").Value = 0

$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 10
$ws.Range("H208").Value = 1
